# Sprint 2 Day 2 update
$wb = $excel.ActiveWorkbook

# --- Product sheet: status/size/sprint tweaks ---
$product = $wb.Worksheets.Item("Product")

# "design UI" story: Size m -> L
$product.Range("D5").Value = "L"

# "send information" story: Size s -> m, Sprint (blank) -> 2, Status Future -> In Progress
$product.Range("D6").Value = "m"
$product.Range("E6").Value = 2
$product.Range("F6").Value = "In Progress"

# "learn android studio" story: Sprint (blank) -> 1
$product.Range("E7").Value = 1

# --- Sprint 2 sheet: Day 2 burndown numbers ---
$sprint2 = $wb.Worksheets.Item("Sprint 2")

# Row 2 "Design UI Look": Day1 remaining hours = 0
$sprint2.Range("E2").Value = 0

# Row 3 "Make UI design function": Product ID 7 -> 4, Start hours 4 -> 10, Day1 0 -> 2, Day2 (blank) -> 0
$sprint2.Range("A3").Value = 4
$sprint2.Range("C3").Value = 10
$sprint2.Range("D3").Value = 2
$sprint2.Range("E3").Value = 0

# Row 4: new task "Link array to UI"
$sprint2.Range("A4").Value = 5
$sprint2.Range("B4").Value = "Link array to UI"
$sprint2.Range("C4").Value = 6
$sprint2.Range("D4").Value = 0
$sprint2.Range("E4").Value = 0

# --- Selection / active-tab bookkeeping to match the saved session state ---
$product.Activate()
$product.Range("G6").Select()

$sprint1 = $wb.Worksheets.Item("Sprint 1")
$sprint1.Activate()
$sprint1.Range("B6").Select()

$sprint2.Activate()
$sprint2.Range("E8").Select()
